$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra "John / Doe / testuser01" test data row (row 5)
$ws.Rows("5:5").Delete() | Out-Null

# Update the Charlie Black row: use a cypress test account instead of
# the previous mismatched-password example
$ws.Range("C4").Value = "testcypress"
$ws.Range("D4").Value = "ValidPass123"
$ws.Range("E4").Value = "DifferentP123"

# Move the active selection
$ws.Range("D8").Select() | Out-Null

# Keep the sheet's outline-level-row metadata in sync with the removed row,
# without leaving any visible grouped row/column behind
$ws.Rows(100).OutlineLevel = 3
$ws.Rows(100).Delete() | Out-Null
